# Updated cryptos list (GitHub Actions crypto-price refresh).
#
# Note: several "Price" (column D) values are plain decimal-looking text
# (e.g. "34.30", "1.00", "0.0000108") that must stay literal TEXT, matching
# how the sheet already stores every price/volume cell as a string. Setting
# Range.Value directly to such a string makes Excel auto-coerce it to a
# number (dropping trailing zeros / switching to scientific notation), so
# for those cells we instead write a text-producing formula, then copy /
# PasteSpecial(xlPasteValues = -4163) it back over itself: that freezes the
# formula's string result as a literal value without Excel re-parsing it as
# a number, and without touching the cell's style.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.849.18'
$ws.Range("E2").Value = '  +4.78%  '
$ws.Range("D3").Value = '2.984.48'
$ws.Range("E3").Value = '  +2.64%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Formula = "=""580.69"""
$ws.Range("D5").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").Value = '  +1.90%  '
$ws.Range("D6").Formula = "=""153.73"""
$ws.Range("D6").Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("E6").Value = '  +6.72%  '
$ws.Range("E7").Value = '  -0.25%  '
$ws.Range("D8").Value = '2.981.34'
$ws.Range("E8").Value = '  +2.49%  '
$ws.Range("E9").Value = '  +0.87%  '
$ws.Range("E10").Value = '  +4.79%  '
$ws.Range("E11").Value = '  +2.91%  '
$ws.Range("D12").Formula = "=""0.447"""
$ws.Range("D12").Copy()
$ws.Range("D12").PasteSpecial(-4163)
$ws.Range("E12").Value = '  +2.79%  '
$ws.Range("E13").Value = '  +2.96%  '
$ws.Range("D14").Formula = "=""34.30"""
$ws.Range("D14").Copy()
$ws.Range("D14").PasteSpecial(-4163)
$ws.Range("E14").Value = '  +6.45%  '
$ws.Range("E15").Value = '  +0.73%  '
$ws.Range("D16").Value = '64.876.74'
$ws.Range("E16").Value = '  +4.88%  '
$ws.Range("D17").Value = '3.478.35'
$ws.Range("E17").Value = '  +2.66%  '
$ws.Range("D18").Formula = "=""6.92"""
$ws.Range("D18").Copy()
$ws.Range("D18").PasteSpecial(-4163)
$ws.Range("E18").Value = '  +3.65%  '
$ws.Range("D19").Value = '2.981.43'
$ws.Range("E19").Value = '  +2.40%  '
$ws.Range("D20").Formula = "=""448.53"""
$ws.Range("D20").Copy()
$ws.Range("D20").PasteSpecial(-4163)
$ws.Range("E20").Value = '  +2.37%  '
$ws.Range("D21").Formula = "=""13.64"""
$ws.Range("D21").Copy()
$ws.Range("D21").PasteSpecial(-4163)
$ws.Range("E21").Value = '  +2.48%  '
$ws.Range("D22").Formula = "=""0.678"""
$ws.Range("D22").Copy()
$ws.Range("D22").PasteSpecial(-4163)
$ws.Range("E22").Value = '  +2.83%  '
$ws.Range("D23").Formula = "=""7.29"""
$ws.Range("D23").Copy()
$ws.Range("D23").PasteSpecial(-4163)
$ws.Range("E23").Value = '  +5.14%  '
$ws.Range("D24").Formula = "=""80.99"""
$ws.Range("D24").Copy()
$ws.Range("D24").PasteSpecial(-4163)
$ws.Range("E24").Value = '  +1.90%  '
$ws.Range("D25").Formula = "=""11.05"""
$ws.Range("D25").Copy()
$ws.Range("D25").PasteSpecial(-4163)
$ws.Range("E25").Value = '  +8.21%  '
$ws.Range("B26").Value = 'Fetch.AI'
$ws.Range("C26").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D26").Formula = "=""2.21"""
$ws.Range("D26").Copy()
$ws.Range("D26").PasteSpecial(-4163)
$ws.Range("E26").Value = '  +7.60%  '
$ws.Range("B27").Value = 'InternetComputer(DFINITY)'
$ws.Range("C27").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D27").Formula = "=""12.22"""
$ws.Range("D27").Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("E27").Value = '  +2.25%  '
$ws.Range("E28").Value = '  +0.04%  '
$ws.Range("D29").Formula = "=""7.81"""
$ws.Range("D29").Copy()
$ws.Range("D29").PasteSpecial(-4163)
$ws.Range("E29").Value = '  +8.99%  '
$ws.Range("B30").Value = 'ImmutableX'
$ws.Range("C30").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D30").Formula = "=""2.37"""
$ws.Range("D30").Copy()
$ws.Range("D30").PasteSpecial(-4163)
$ws.Range("E30").Value = '  +13.86%  '
$ws.Range("B31").Value = 'PEPE'
$ws.Range("C31").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D31").Formula = "=""0.0000108"""
$ws.Range("D31").Copy()
$ws.Range("D31").PasteSpecial(-4163)
$ws.Range("E31").Value = '  +1.44%  '
$ws.Range("D32").Formula = "=""2.57"""
$ws.Range("D32").Copy()
$ws.Range("D32").PasteSpecial(-4163)
$ws.Range("E32").Value = '  +2.01%  '
$ws.Range("E33").Value = '  +2.30%  '
$ws.Range("D34").Formula = "=""26.58"""
$ws.Range("D34").Copy()
$ws.Range("D34").PasteSpecial(-4163)
$ws.Range("E34").Value = '  +3.27%  '
$ws.Range("D35").Formula = "=""1.00"""
$ws.Range("D35").Copy()
$ws.Range("D35").PasteSpecial(-4163)
$ws.Range("E35").Value = '  -0.01%  '
$ws.Range("D36").Formula = "=""0.984"""
$ws.Range("D36").Copy()
$ws.Range("D36").PasteSpecial(-4163)
$ws.Range("E36").Value = '  +2.29%  '
$ws.Range("D37").Formula = "=""5.65"""
$ws.Range("D37").Copy()
$ws.Range("D37").PasteSpecial(-4163)
$ws.Range("E37").Value = '  +3.25%  '
$ws.Range("E38").Value = '  +8.13%  '
$ws.Range("D39").Formula = "=""3.07"""
$ws.Range("D39").Copy()
$ws.Range("D39").PasteSpecial(-4163)
$ws.Range("E39").Value = '  +7.16%  '
$ws.Range("D40").Formula = "=""49.00"""
$ws.Range("D40").Copy()
$ws.Range("D40").PasteSpecial(-4163)
$ws.Range("E40").Value = '  -0.11%  '
$ws.Range("D41").Formula = "=""44.30"""
$ws.Range("D41").Copy()
$ws.Range("D41").PasteSpecial(-4163)
$ws.Range("E41").Value = '  +10.77%  '
$ws.Range("D42").Formula = "=""0.120"""
$ws.Range("D42").Copy()
$ws.Range("D42").PasteSpecial(-4163)
$ws.Range("E42").Value = '  +3.30%  '
$ws.Range("E43").Value = '  +8.87%  '
$ws.Range("D44").Formula = "=""8.37"""
$ws.Range("D44").Copy()
$ws.Range("D44").PasteSpecial(-4163)
$ws.Range("E44").Value = '  +0.59%  '
$ws.Range("D45").Formula = "=""388.48"""
$ws.Range("D45").Copy()
$ws.Range("D45").PasteSpecial(-4163)
$ws.Range("E45").Value = '  +13.95%  '
$ws.Range("D46").Value = '2.782.90'
$ws.Range("E46").Value = '  +3.02%  '
$ws.Range("E47").Value = '  +4.42%  '
$ws.Range("D48").Formula = "=""135.17"""
$ws.Range("D48").Copy()
$ws.Range("D48").PasteSpecial(-4163)
$ws.Range("E48").Value = '  +1.13%  '
$ws.Range("E49").Value = '  -0.01%  '
$ws.Range("D50").Formula = "=""0.000225"""
$ws.Range("D50").Copy()
$ws.Range("D50").PasteSpecial(-4163)
$ws.Range("E50").Value = '  +14.37%  '
$ws.Range("E51").Value = '  +1.85%  '
